# AGGP247 Schedule S20 - "Updated Sylubus and Schedule"
#
# The shared-strings table was reshuffled (Excel does this bookkeeping
# automatically as cell text changes), so the only edits that matter are
# the actual visible content/format changes below:
#   - K7:  "Circle 1"  -> "Circles"
#   - M7:  "Lab 2 DUE" -> "No Class 17th, Lab 2 DUE"
#   - K12: "Gravity  & Forces" -> "Gravity & Forces" (extra space removed)
#   - K13: "Gravity  & Forces" -> "Gravity & Forces" (extra space removed)
#   - M16: (blank) -> "Project Review" (and now bold, like the rest of col M)
#   - M17: "Final Due" -> (blank)
#   - M18: (new cell) -> "Final Due" (bold, like the rest of col M)
#   - M6:  new (blank) bold-style cell added
#   - B7:  highlighted (red fill) like the Spring-Break week row
#   - Selection moved from K13 to K25, and the frozen/scrolled
#     left column (topLeftCell=C1) is cleared back to normal.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text edits -----------------------------------------------------
$ws.Range("K7").Value = "Circles"
$ws.Range("M7").Value = "No Class 17th, Lab 2 DUE"

$ws.Range("K12").Value = "Gravity & Forces"
$ws.Range("K13").Value = "Gravity & Forces"

$ws.Range("M16").Value = "Project Review"
$ws.Range("M17").Value = ""
$ws.Range("M18").Value = "Final Due"

# --- Formatting edits -------------------------------------------------
# Give the new/changed Notes-column cells the same bold style already
# used throughout column M (copy format from M19, which already carries
# that style and no value).
$ws.Range("M19").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("M18").PasteSpecial(-4122)

# Highlight B7 (Feb 17 week) the same way B11 (Spring Break week) is
# already highlighted.
$ws.Range("B11").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View/selection edits ---------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
$ws.Range("K25").Select() | Out-Null
